# Update the "Main Info" sheet:
#  - B16 (Main Link -> 4G+Cellular) flips from TRUE to FALSE
#  - B24 (Backup Link -> 4G+Cellular) flips from FALSE to TRUE
#  - The active selection moves to G16 (reflecting that the edit happened there)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

$ws.Range("B16").Value = $false
$ws.Range("B24").Value = $true

$ws.Activate()
[void]$ws.Range("G16").Select()
